$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four trailing listing rows (rows 3-6); only the header row
# and the first data row remain.
$ws.Rows("3:6").Delete()

# Drop every existing hyperlink (and its backing relationship) so we can
# rebuild just the single one the surviving row needs.
$ws.Hyperlinks.Delete()

# C2:E2 hold values that look numeric/date-like ("2021-11-26", "1788",
# "60067867890") but must stay text (shared strings), matching the
# source data. Force text interpretation via NumberFormat, assign, then
# restore the Normal style so no stray number format sticks to the cells.
$ws.Range("C2:E2").NumberFormat = "@"

$ws.Range("A2").Value = "Особняк на Александрова, 18"
$ws.Range("B2").Value = "https://osobnyaki.com//na-nikoloyamskoy-49s2"
$ws.Range("C2").Value = "2021-11-26"
$ws.Range("D2").Value = "1788"
$ws.Range("E2").Value = "60067867890"
$ws.Range("F2").Value = 54152

$ws.Range("C2:E2").Style = "Normal"

# Re-create the hyperlink on B2 pointing at the new address, then restore
# the plain built-in Hyperlink cell style (Hyperlinks.Add nudges the font
# flag on the applied style, which would otherwise mint a redundant xf).
$ws.Hyperlinks.Add($ws.Range("B2"), "https://osobnyaki.com//na-nikoloyamskoy-49s2")
$ws.Range("B2").Style = "Hyperlink"
